$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 138, shifting existing rows
# 138-145 down to 139-146 (values + formatting move with the rows).
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with the new weekly record.
$ws.Range("A138").Value = 3
$ws.Range("B138").Value = "Femacal de La Calera"
$ws.Range("C138").Value = "Coquimbo"
$ws.Range("D138").Value = 44610
$ws.Range("E138").Value = 5
$ws.Range("F138").Value = 100112030
$ws.Range("G138").Value = "Poroto granado"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 70
$ws.Range("K138").Value = 24000
$ws.Range("L138").Value = 25000
$ws.Range("M138").Value = 24500
$ws.Range("N138").Value = "$/malla 25 kilos"
$ws.Range("O138").Value = "Provincia de Petorca"
$ws.Range("P138").Value = 980
$ws.Range("Q138").Value = 25
$ws.Range("R138").Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Range("D138").NumberFormat = $ws.Range("D139").NumberFormat
